$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New quarter row appended below the existing data (row 63).
# Build the date label as a formula first, then convert it in-place to a
# plain value via copy/paste-special so Excel doesn't auto-coerce the
# "01-04-2021" literal into a date serial (which would also pull in an
# extra number-format style not present in the target workbook).
$ws.Range("A63").Formula = "=""01-04-2021"""
$ws.Range("A63").Copy()
$ws.Range("A63").PasteSpecial(-4163)

$ws.Range("B63").Value = 111408472
$ws.Range("C63").Value = 27525294
